$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dist_coûts")

# Rows 2-6, column Y: 75 -> 2.25
$ws.Range("Y2:Y6").Value = 2.25

# Row 8, column Y: 0 -> 1
$ws.Range("Y8").Value = 1

# Row 10, columns Y:AA: 2000/2100/2200 -> 0
$ws.Range("Y10:AA10").Value = 0

# Move the active selection to Y10 (matches the saved sheet view state)
$ws.Range("Y10").Select()
